{"js": "// The title paragraph text changes from \u00ab\u0420\u0430\u0431\u043e\u0442\u0430 \u0432 \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u043e\u043c \u0440\u0435\u0436\u0438\u043c\u0435\u00bb\n// to \u00ab\u0418\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u043e\u0432 \u0441\u043e\u0440\u0442\u0438\u0440\u043e\u0432\u043a\u0438\u00bb, and the `_GoBack` bookmark\n// (previously sitting right after the \"\u21165\" run) moves down to sit right\n// before the \"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u2026\" run in the paragraph that follows the title.\n\nconst body = context.document.body;\n\n// 1) Replace the title text, preserving the run/paragraph formatting that\n//    is already on the matched range.\nconst titleHits = body.search(\"\u00ab\u0420\u0430\u0431\u043e\u0442\u0430 \u0432 \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u043e\u043c \u0440\u0435\u0436\u0438\u043c\u0435\u00bb\", { matchCase: false });\ntitleHits.load(\"text\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"\u00ab\u0418\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u043e\u0432 \u0441\u043e\u0440\u0442\u0438\u0440\u043e\u0432\u043a\u0438\u00bb\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Move the `_GoBack` bookmark from the end of the \"\u2026\u0420\u0410\u0411\u041e\u0422\u0415 \u21165\" paragraph\n//    to the very start of the following \"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u2026\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst disciplineHits = body.search(\"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u00ab\u041e\u0441\u043d\u043e\u0432\u044b \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u0438\u0437\u0430\u0446\u0438\u0438 \u0438 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f\u00bb\", { matchCase: false });\ndisciplineHits.load(\"text\");\nawait context.sync();\n\nif (disciplineHits.items.length > 0) {\n  const startOfDiscipline = disciplineHits.items[0].getRange(\"Start\");\n  startOfDiscipline.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The title paragraph text changes from \u00ab\u0420\u0430\u0431\u043e\u0442\u0430 \u0432 \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u043e\u043c \u0440\u0435\u0436\u0438\u043c\u0435\u00bb\n# to \u00ab\u0418\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u043e\u0432 \u0441\u043e\u0440\u0442\u0438\u0440\u043e\u0432\u043a\u0438\u00bb, and the `_GoBack` bookmark\n# (previously sitting right after the \"\u21165\" run) moves down to sit right\n# before the \"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u2026\" run in the paragraph that follows the title.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the title text via Find/Replace - formatting of the matched\n#    range (bold, Times New Roman, 14pt) carries over to the new text.\n$find = $d.Content.Find\n$find.Text = \"\u00ab\u0420\u0430\u0431\u043e\u0442\u0430 \u0432 \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u043e\u043c \u0440\u0435\u0436\u0438\u043c\u0435\u00bb\"\n$find.Replacement.Text = \"\u00ab\u0418\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u0435 \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u043e\u0432 \u0441\u043e\u0440\u0442\u0438\u0440\u043e\u0432\u043a\u0438\u00bb\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Move the `_GoBack` bookmark from the end of the \"\u2026\u0420\u0410\u0411\u041e\u0422\u0415 \u21165\" paragraph\n#    to the very start of the following \"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u2026\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$target = $d.Content\n$targetFind = $target.Find\n$targetFind.Text = \"\u041f\u041e \u0414\u0418\u0421\u0426\u0418\u041f\u041b\u0418\u041d\u0415 \u00ab\u041e\u0441\u043d\u043e\u0432\u044b \u0430\u043b\u0433\u043e\u0440\u0438\u0442\u043c\u0438\u0437\u0430\u0446\u0438\u0438 \u0438 \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f\u00bb\"\n$targetFind.Execute() | Out-Null\n$target.Collapse(1) | Out-Null\n\n$d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n"}
